$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.647.24'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').Value = '3.775.33'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'439.33"
$ws.Range('E5').Value = '  +4.55%  '
$ws.Range('D6').Value = "'142.93"
$ws.Range('E6').Value = '  +10.79%  '
$ws.Range('D7').Value = "'0.620"
$ws.Range('E7').Value = '  +3.31%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').Value = "'0.731"
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('D10').Value = "'0.151"
$ws.Range('E10').Value = '  -7.39%  '
$ws.Range('D11').Value = "'0.0000309"
$ws.Range('E11').Value = '  -9.23%  '
$ws.Range('D12').Value = "'42.71"
$ws.Range('E12').Value = '  +6.39%  '
$ws.Range('D13').Value = "'10.40"
$ws.Range('E13').Value = '  +5.01%  '
$ws.Range('D14').Value = '4.393.33'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').Value = "'14.78"
$ws.Range('E15').Value = '  -7.63%  '
$ws.Range('D16').Value = '3.823.07'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = "'19.83"
$ws.Range('E18').Value = '  +2.31%  '
$ws.Range('E19').Value = '  +6.53%  '
$ws.Range('D20').Value = '66.725.06'
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').Value = "'410.82"
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('D22').Value = "'14.46"
$ws.Range('E22').Value = '  +1.28%  '
$ws.Range('D23').Value = "'3.26"
$ws.Range('E23').Value = '  +9.03%  '
$ws.Range('D24').Value = "'85.24"
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('D25').Value = "'3.40"
$ws.Range('E25').Value = '  +6.29%  '
$ws.Range('D26').Value = "'36.79"
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = "'5.61"
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').Value = "'9.65"
$ws.Range('E28').Value = '  +33.22%  '
$ws.Range('D29').Value = "'9.70"
$ws.Range('E29').Value = '  +3.35%  '
$ws.Range('D30').Value = "'730.57"
$ws.Range('E30').Value = '  +5.58%  '
$ws.Range('D31').Value = "'13.82"
$ws.Range('E31').Value = '  +13.03%  '
$ws.Range('E32').Value = '  +11.70%  '
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('D34').Value = "'42.58"
$ws.Range('E34').Value = '  +12.56%  '
$ws.Range('D35').Value = "'0.157"
$ws.Range('E35').Value = '  +4.58%  '
$ws.Range('E36').Value = '  +26.03%  '
$ws.Range('D37').Value = "'56.18"
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('D38').Value = "'0.999"
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('D39').Value = "'0.0476"
$ws.Range('E39').Value = '  +6.04%  '
$ws.Range('D40').Value = "'2.71"
$ws.Range('E40').Value = '  +37.69%  '
$ws.Range('D41').Value = "'2.89"
$ws.Range('E41').Value = '  -1.36%  '
$ws.Range('D42').Value = "'3.34"
$ws.Range('E42').Value = '  +7.77%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = "'0.139"
$ws.Range('E44').Value = '  +3.98%  '
$ws.Range('B45').Value = 'PEPE'
$ws.Range('C45').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D45').Value = '0.0₃0664'
$ws.Range('E45').Value = '  -11.70%  '
$ws.Range('D46').Value = "'0.329"
$ws.Range('E46').Value = '  +15.18%  '
$ws.Range('D47').Value = "'3.29"
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = "'2.66"
$ws.Range('E48').Value = '  +5.20%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = "'2.08"
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').Value = "'143.07"
$ws.Range('E50').Value = '  -0.59%  '
$ws.Range('D51').Value = "'2.81"
$ws.Range('E51').Value = '  +2.78%  '

foreach ($addr in @('D5','D6','D7','D9','D10','D11','D12','D13','D15','D18','D21','D22','D23','D24','D25','D26','D27','D28','D29','D30','D31','D34','D35','D37','D38','D39','D40','D41','D42','D43','D44','D46','D47','D48','D49','D50','D51')) {
    $ws.Range($addr).Style = "Normal"
}
